$d = $word.ActiveDocument

$replacements = @(
    @("483×3=", "616×8="),
    @("213×9=", "563×8="),
    @("404×9=", "203×9="),
    @("232×9=", "710×9="),
    @("292×4=", "614×2="),
    @("145×7=", "260×2="),
    @("252×6=", "455×2="),
    @("147×6=", "685×9="),
    @("726×8=", "158×9="),
    @("119×9=", "808×2="),
    @("533×5=", "257×2="),
    @("503×7=", "249×3="),
    @("933×7=", "479×9="),
    @("757×7=", "732×4="),
    @("901×3=", "116×5="),
    @("448×8=", "693×3="),
    @("521×5=", "541×8="),
    @("978×7=", "188×3="),
    @("616×3=", "432×8="),
    @("249×8=", "400×6="),
    @("202×8=", "634×8="),
    @("678×3=", "687×2="),
    @("118×3=", "424×4="),
    @("282×4=", "140×9="),
    @("245×6=", "222×7=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
